$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Orange"
$ws.Range("A7").Value = "Orange"
$ws.Range("A8").Value = "Orange"

$ws.Range("A8").Select()
